$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the IG "Date" property ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-06-04T08:55:54+00:00"

# --- Elements sheet: Extension.value[x] row (row 6) ---
$elements = $wb.Worksheets.Item("Elements")

# Type(s): code -> CodeableConcept (trailing newline, matches the existing
# "Type(s)" column convention in this sheet)
$elements.Range("K6").Value = "CodeableConcept`n"

# New required binding + value set for the value[x] element
$elements.Range("X6").Value = "required"
$elements.Range("Y6").Value = ""
$elements.Range("Z6").Value = "http://ltsi.univ-rennes.fr/ValueSet/siph-statutcomposantadm-oncofair-valueset"

# Widen the Type(s) / Binding Value Set columns to fit the new content
$elements.Columns.Item(11).ColumnWidth = 16.25
$elements.Columns.Item(26).ColumnWidth = 71.92

# Keep row 6 at its normal auto height (avoid a stray custom row height
# from the longer wrapped text we just entered)
$elements.Rows.Item(6).AutoFit()
